$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D (shifts existing D:K data to F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Copy number formats/styles from the (now-shifted) old D:E columns (F:G) onto the new blank D:E columns.
# Done per contiguous block so we don't materialize cells on the blank spacer rows (36, 78).
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F37:G77").Copy()
$ws.Range("D37:E77").PasteSpecial(-4122)
$ws.Range("F79:G102").Copy()
$ws.Range("D79:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D, E) with the latest reported figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 443200
$ws.Range("E8").Value = 102300
$ws.Range("D9").Value = 35100
$ws.Range("E9").Value = 31700
$ws.Range("D10").Value = 408100
$ws.Range("E10").Value = 70600
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 102000
$ws.Range("E14").Value = 9500
$ws.Range("D15").Value = 48300
$ws.Range("E15").Value = 30300
$ws.Range("D17").Value = 190400
$ws.Range("E17").Value = 76200
$ws.Range("D18").Value = 252800
$ws.Range("E18").Value = 26100
$ws.Range("D20").Value = -14500
$ws.Range("E20").Value = 13300
$ws.Range("D21").Value = 286600
$ws.Range("E21").Value = 69700
$ws.Range("D22").Value = 20100
$ws.Range("E22").Value = 20400
$ws.Range("D23").Value = 218200
$ws.Range("E23").Value = 19000
$ws.Range("D24").Value = -100
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 218300
$ws.Range("E26").Value = 19000
$ws.Range("D27").Value = 218300
$ws.Range("E27").Value = 19000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 14500
$ws.Range("E32").Value = -13300
$ws.Range("D33").Value = 218300
$ws.Range("E33").Value = 19000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 218300
$ws.Range("E35").Value = 19000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 2400
$ws.Range("E41").Value = 113000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 97600
$ws.Range("E43").Value = 91300
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 128500
$ws.Range("E45").Value = 20500
$ws.Range("D46").Value = 228400
$ws.Range("E46").Value = 224700
$ws.Range("D47").Value = 5100
$ws.Range("E47").Value = 5200
$ws.Range("D48").Value = 1202700
$ws.Range("E48").Value = 818000
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 67400
$ws.Range("E52").Value = 21000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1503600
$ws.Range("E54").Value = 1068900
$ws.Range("D57").Value = 135500
$ws.Range("E57").Value = 122300
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("D59").Value = 96000
$ws.Range("E59").Value = 93100
$ws.Range("D60").Value = 231500
$ws.Range("E60").Value = 215400
$ws.Range("D61").Value = 830200
$ws.Range("E61").Value = 789500
$ws.Range("D62").Value = 12100
$ws.Range("E62").Value = 52800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1073800
$ws.Range("E66").Value = 1057700
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -796900
$ws.Range("E72").Value = -1015200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 429900
$ws.Range("E76").Value = 11200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 218300
$ws.Range("E81").Value = 19000
$ws.Range("D83").Value = 48300
$ws.Range("E83").Value = 30300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 117800
$ws.Range("E89").Value = 62800
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -164400
$ws.Range("E94").Value = -150400
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -64000
$ws.Range("E100").Value = -300
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -110600
$ws.Range("E102").Value = -88000

# Restated figures for a handful of historical quarters that changed in this refresh
$ws.Range("J20").Value = -100
$ws.Range("J21").Value = 43900
$ws.Range("J22").Value = 16400
$ws.Range("J32").Value = 100
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 42700
